$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text, preserving existing (unstyled) format
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" '27.088.66'
Set-TextCell $ws "E2" '  +0.58%  '

# Row 3
Set-TextCell $ws "D3" '1.890.41'
Set-TextCell $ws "E3" '  +1.47%  '

# Row 4
Set-TextCell $ws "D4" '1.0000'
Set-TextCell $ws "E4" '  +0.02%  '

# Row 5
Set-TextCell $ws "D5" '307.63'
Set-TextCell $ws "E5" '  +0.89%  '

# Row 6
Set-TextCell $ws "D6" '1.000'
Set-TextCell $ws "E6" '  +0.03%  '

# Row 7
Set-TextCell $ws "D7" '0.5143'
Set-TextCell $ws "E7" '  +1.44%  '

# Row 8
Set-TextCell $ws "D8" '0.3742'
Set-TextCell $ws "E8" '  +3.13%  '

# Row 9
Set-TextCell $ws "D9" '0.07212'
Set-TextCell $ws "E9" '  +0.54%  '

# Row 10
Set-TextCell $ws "B10" 'Solana'
Set-TextCell $ws "C10" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws "D10" '21.16'
Set-TextCell $ws "E10" '  +2.16%  '

# Row 11
Set-TextCell $ws "B11" 'Polygon'
Set-TextCell $ws "C11" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws "D11" '0.9052'
Set-TextCell $ws "E11" '  +0.93%  '

# Row 12
Set-TextCell $ws "D12" '0.07639'
Set-TextCell $ws "E12" '  +2.57%  '

# Row 13
Set-TextCell $ws "D13" '1.886.20'
Set-TextCell $ws "E13" '  +1.25%  '

# Row 14
Set-TextCell $ws "D14" '94.76'
Set-TextCell $ws "E14" '  +2.25%  '

# Row 15
Set-TextCell $ws "D15" '5.273'
Set-TextCell $ws "E15" '  +0.61%  '

# Row 16
Set-TextCell $ws "D16" '1.001'
Set-TextCell $ws "E16" '  +0.04%  '

# Row 17
Set-TextCell $ws "D17" '0.000008497'
Set-TextCell $ws "E17" '  +0.08%  '

# Row 18
Set-TextCell $ws "D18" '14.46'
Set-TextCell $ws "E18" '  +2.25%  '

# Row 19
Set-TextCell $ws "E19" '  -0.01%  '

# Row 20
Set-TextCell $ws "D20" '27.126.13'
Set-TextCell $ws "E20" '  +0.60%  '

# Row 21
Set-TextCell $ws "D21" '5.074'
Set-TextCell $ws "E21" '  +0.95%  '

# Row 22
Set-TextCell $ws "D22" '2.119.86'
Set-TextCell $ws "E22" '  +1.41%  '

# Row 23
Set-TextCell $ws "E23" '  +2.16%  '

# Row 24
Set-TextCell $ws "D24" '6.415'
Set-TextCell $ws "E24" '  -0.27%  '

# Row 25
Set-TextCell $ws "D25" '146.06'
Set-TextCell $ws "E25" '  -1.30%  '

# Row 26
Set-TextCell $ws "D26" '1.788'
Set-TextCell $ws "E26" '  -0.36%  '

# Row 27
Set-TextCell $ws "D27" '2.216'
Set-TextCell $ws "E27" '  +7.37%  '

# Row 28
Set-TextCell $ws "E28" '  +1.18%  '

# Row 29
Set-TextCell $ws "E29" '  +1.14%  '

# Row 30
Set-TextCell $ws "D30" '4.859'
Set-TextCell $ws "E30" '  +4.05%  '

# Row 31
Set-TextCell $ws "D31" '4.960'
Set-TextCell $ws "E31" '  +5.99%  '

# Row 32
Set-TextCell $ws "D32" '0.09195'
Set-TextCell $ws "E32" '  -0.67%  '

# Row 33
Set-TextCell $ws "D33" '0.05086'
Set-TextCell $ws "E33" '  +0.02%  '

# Row 34
Set-TextCell $ws "D34" '1.235'
Set-TextCell $ws "E34" '  +7.04%  '

# Row 35
Set-TextCell $ws "D35" '0.7686'
Set-TextCell $ws "E35" '  +2.72%  '

# Row 36
Set-TextCell $ws "D36" '2.977'
Set-TextCell $ws "E36" '  -0.45%  '

# Row 37
Set-TextCell $ws "D37" '3.290'
Set-TextCell $ws "E37" '  +0.58%  '

# Row 38
Set-TextCell $ws "D38" '2.606'
Set-TextCell $ws "E38" '  +3.03%  '

# Row 39
Set-TextCell $ws "D39" '0.02000'
Set-TextCell $ws "E39" '  -0.10%  '

# Row 40
Set-TextCell $ws "D40" '0.5603'
Set-TextCell $ws "E40" '  +3.18%  '

# Row 41
Set-TextCell $ws "D41" '1.075'
Set-TextCell $ws "E41" '  -0.58%  '

# Row 42
Set-TextCell $ws "D42" '6.652'
Set-TextCell $ws "E42" '  +2.37%  '

# Row 43
Set-TextCell $ws "D43" '8.956'
Set-TextCell $ws "E43" '  +4.47%  '

# Row 44
Set-TextCell $ws "D44" '117.99'
Set-TextCell $ws "E44" '  +0.00%  '

# Row 45
Set-TextCell $ws "E45" '  +3.08%  '

# Row 46
Set-TextCell $ws "D46" '0.4804'
Set-TextCell $ws "E46" '  +2.74%  '

# Row 47
Set-TextCell $ws "D47" '10.21'
Set-TextCell $ws "E47" '  +1.61%  '

# Row 48
Set-TextCell $ws "D48" '0.9995'
Set-TextCell $ws "E48" '  +0.00%  '

# Row 49
Set-TextCell $ws "D49" '1.592'
Set-TextCell $ws "E49" '  +1.70%  '

# Row 50
Set-TextCell $ws "D50" '37.51'
Set-TextCell $ws "E50" '  +1.60%  '

# Row 51
Set-TextCell $ws "D51" '63.94'
Set-TextCell $ws "E51" '  +1.37%  '
